# Apply text updates described in the commit:
# "From 1.2.4 to 1.2.5 change and minor updates"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Version: 1.0 -> 1.2.5  (D2)
$ws.Range("D2").Value = "1.2.5"

# 2) Precondition text fix (typo + period) - appears once per test case (TC1..TC5)
$newPrecondition = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B8").Value  = $newPrecondition
$ws.Range("B15").Value = $newPrecondition
$ws.Range("B22").Value = $newPrecondition
$ws.Range("B29").Value = $newPrecondition
$ws.Range("B36").Value = $newPrecondition

# 3) Expected results - add trailing period
$ws.Range("D10").Value = "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária."
$ws.Range("D17").Value = "SYSTEM Apresenta a tela de Analisar Prestação de Contas."
$ws.Range("D24").Value = "SYSTEM Apresenta a tela de Detalhar Diárias."

# 4) TC5 step text rewording
$ws.Range("B38").Value = "Beneficiário Acessa o caso de uso através do menu."
